# Esecizio migliorato con BETWEEN.
# Replace the `LIKE '1975-%';` clause with `BETWEEN '1975-01-01' AND '1975-31-12';`

$d = $word.ActiveDocument

$rng = $d.Content.Duplicate
$found = $rng.Find.Execute("LIKE '1975-%';", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Assigning .Text directly (instead of using Find's Replace argument) avoids
    # the "smart quotes" autocorrect substitution and collapses the found range
    # into a single run using the formatting of the first character.
    $rng.Text = "BETWEEN '1975-01-01' AND '1975-31-12';"
}
